$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 143, shifting existing rows 143:162 down to 144:163
$ws.Rows.Item(143).Insert()

# Populate the new row 143 with the new record
$ws.Range("A143").Value = 11
$ws.Range("B143").Value = "Vega Monumental Concepción"
$ws.Range("C143").Value = "Bíobío"
$ws.Range("D143").Value = 44505
$ws.Range("E143").Value = 8
$ws.Range("F143").Value = 100112017
$ws.Range("G143").Value = "Apio"
$ws.Range("H143").Value = "Americana (o)"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 220
$ws.Range("K143").Value = 7500
$ws.Range("L143").Value = 8000
$ws.Range("M143").Value = 7773
$ws.Range("N143").Value = "$/docena de matas"
$ws.Range("O143").Value = "Provincia del Elquí"
$ws.Range("P143").Value = 1296
$ws.Range("Q143").Value = 6
$ws.Range("R143").Value = "Hortaliza"

# Match the date-number-format style used by column D elsewhere
$ws.Range("D143").NumberFormat = $ws.Range("D142").NumberFormat
